# Trading update: 2026-02-17 13:02:13
# Append a new trade row (row 4) to both the "All Trades" sheet and the
# "MarketMaking" sheet.

$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "MarketMaking")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item(4, 1).Value = 3

    # Force Date/Time columns to remain plain text instead of being
    # auto-converted into date/time serial numbers.
    $ws.Cells.Item(4, 2).NumberFormat = "@"
    $ws.Cells.Item(4, 2).Value = "2026-02-17"
    $ws.Cells.Item(4, 3).NumberFormat = "@"
    $ws.Cells.Item(4, 3).Value = "13:02:05"

    $ws.Cells.Item(4, 4).Value = "MarketMaking"
    $ws.Cells.Item(4, 5).Value = "UP"
    $ws.Cells.Item(4, 6).Value = 0.54
    # G4 (Exit Price) stays empty - trade is still OPEN.
    $ws.Cells.Item(4, 8).Value = "OPEN"
    $ws.Cells.Item(4, 9).Value = 0
    $ws.Cells.Item(4, 10).Value = 0
    $ws.Cells.Item(4, 11).Value = 99.86
    $ws.Cells.Item(4, 12).Value = 0
    $ws.Cells.Item(4, 13).Value = 0
    $ws.Cells.Item(4, 14).Value = 0.6
    $ws.Cells.Item(4, 15).Value = "Normal spread capture: 19600 bps"
    # P4 (Exit Reason) stays empty - trade is still OPEN.
    $ws.Cells.Item(4, 17).Value = 0
}
